$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 431.7736
$ws.Range("J17").Value = 438.90195
$ws.Range("L17").Value = 1316.70585
$ws.Range("N17").Value = -1652.70585
$ws.Range("H40").Value = 2741.25
$ws.Range("I40").Value = 2722
$ws.Range("K40").Value = 2722
$ws.Range("M40").Value = -2547
$ws.Range("H86").Value = 9092922
$ws.Range("I86").Value = 14287455
$ws.Range("J86").Value = 2488.75
$ws.Range("K86").Value = 14287455
$ws.Range("L86").Value = 2488.75
$ws.Range("M86").Value = -14286332
$ws.Range("N86").Value = -4734.75
$ws.Range("H89").Value = 9092922
$ws.Range("I89").Value = 14287455
$ws.Range("J89").Value = 2488.75
$ws.Range("K89").Value = 71437275
$ws.Range("L89").Value = 12443.75
$ws.Range("M89").Value = -71431659
$ws.Range("N89").Value = -23675.75
$ws.Range("H92").Value = 1026.5264
$ws.Range("I92").Value = 819.6923
$ws.Range("K92").Value = 819.6923
$ws.Range("M92").Value = 428.3077
$ws.Range("H132").Value = 1557.1025
$ws.Range("I132").Value = 1558.6052
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 4675.8156
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -2145.8156
$ws.Range("N132").Value = -9560
$ws.Range("H137").Value = 286425.2
$ws.Range("I137").Value = 1504.7576
$ws.Range("J137").Value = 808779.25
$ws.Range("K137").Value = 4514.2728
$ws.Range("L137").Value = 2426337.75
$ws.Range("M137").Value = -1964.2728
$ws.Range("N137").Value = -2431437.75
$ws.Range("H141").Value = 2953.9092
$ws.Range("I141").Value = 2249.1
$ws.Range("K141").Value = 6747.299999999999
$ws.Range("M141").Value = -1567.299999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3548
$ws.Range("I2").Value = 2000
$ws.Range("J2").Value = 5096
$ws.Range("K2").Value = 2000
$ws.Range("L2").Value = 5096
$ws.Range("M2").Value = -1887
$ws.Range("N2").Value = -5322
$ws.Range("H32").Value = 4267.493
$ws.Range("I32").Value = 1829.3383
$ws.Range("K32").Value = 1829.3383
$ws.Range("M32").Value = -1542.3383
$ws.Range("H38").Value = 21591.8
$ws.Range("I38").Value = 21591.8
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 21591.8
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -21124.8
$ws.Range("N38").ClearContents()
$ws.Range("H63").Value = 3705.8
$ws.Range("J63").Value = 4531.6
$ws.Range("L63").Value = 4531.6
$ws.Range("N63").Value = -5903.6
$ws.Range("H66").Value = 3705.8
$ws.Range("J66").Value = 4531.6
$ws.Range("L66").Value = 22658
$ws.Range("N66").Value = -29522
$ws.Range("H116").Value = 3548
$ws.Range("I116").Value = 2000
$ws.Range("J116").Value = 5096
$ws.Range("K116").Value = 2000
$ws.Range("L116").Value = 5096
$ws.Range("M116").Value = 294
$ws.Range("N116").Value = -9684
$ws.Range("H132").Value = 3518.35
$ws.Range("I132").Value = 3345.2354
$ws.Range("K132").Value = 10035.7062
$ws.Range("M132").Value = -7505.706200000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5096
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 5096
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 5096
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -5324
$ws.Range("H19").Value = 2000
$ws.Range("J19").Value = 2000
$ws.Range("L19").Value = 2000
$ws.Range("N19").Value = -2346
$ws.Range("H105").Value = 252397.75
$ws.Range("I105").Value = 335464.34
$ws.Range("K105").Value = 335464.34
$ws.Range("M105").Value = -333717.34

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1164.2609
$ws.Range("I22").Value = 878.35297
$ws.Range("J22").Value = 1974.3334
$ws.Range("K22").Value = 878.35297
$ws.Range("L22").Value = 1974.3334
$ws.Range("M22").Value = -528.35297
$ws.Range("N22").Value = -2674.3334
$ws.Range("H32").Value = 1644.2222
$ws.Range("I32").Value = 1828.2858
$ws.Range("J32").Value = 1000
$ws.Range("K32").Value = 1828.2858
$ws.Range("L32").Value = 1000
$ws.Range("M32").Value = -1512.2858
$ws.Range("N32").Value = -1632
$ws.Range("H35").Value = 6417.6
$ws.Range("I35").Value = 3515
$ws.Range("J35").Value = 18028
$ws.Range("K35").Value = 3515
$ws.Range("L35").Value = 18028
$ws.Range("M35").Value = -3221
$ws.Range("N35").Value = -18616
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H68").Value = 28479
$ws.Range("I68").Value = 28479
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 28479
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -27730
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 28479
$ws.Range("I71").Value = 28479
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 85437
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -81693
$ws.Range("N71").ClearContents()
$ws.Range("H74").Value = 60000
$ws.Range("I74").Value = 20000
$ws.Range("J74").Value = 100000
$ws.Range("K74").Value = 20000
$ws.Range("L74").Value = 100000
$ws.Range("M74").Value = -19126
$ws.Range("N74").Value = -101748
$ws.Range("H77").Value = 60000
$ws.Range("I77").Value = 20000
$ws.Range("J77").Value = 100000
$ws.Range("K77").Value = 60000
$ws.Range("L77").Value = 300000
$ws.Range("M77").Value = -55632
$ws.Range("N77").Value = -308736
$ws.Range("H105").Value = 39974.066
$ws.Range("I105").Value = 57068.9
$ws.Range("K105").Value = 57068.9
$ws.Range("M105").Value = -55321.9
$ws.Range("H132").Value = 1247979.1
$ws.Range("I132").Value = 1086728.4
$ws.Range("J132").Value = 2602485.5
$ws.Range("K132").Value = 3260185.2
$ws.Range("L132").Value = 7807456.5
$ws.Range("M132").Value = -3257655.2
$ws.Range("N132").Value = -7812516.5
$ws.Range("H134").Value = 3762307
$ws.Range("I134").Value = 4204549
$ws.Range("J134").Value = 3250
$ws.Range("K134").Value = 12613647
$ws.Range("L134").Value = 9750
$ws.Range("M134").Value = -12611112
$ws.Range("N134").Value = -14820

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 51.285713
$ws.Range("I17").Value = 58.166668
$ws.Range("J17").Value = 10
$ws.Range("K17").Value = 174.500004
$ws.Range("L17").Value = 30
$ws.Range("M17").Value = -5.50000399999999
$ws.Range("N17").Value = -368

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 435543.56
$ws.Range("I10").Value = 7433.3335
$ws.Range("J10").Value = 756626.25
$ws.Range("K10").Value = 7433.3335
$ws.Range("L10").Value = 756626.25
$ws.Range("M10").Value = -7264.3335
$ws.Range("N10").Value = -756964.25
$ws.Range("H40").Value = 14053.333
$ws.Range("I40").Value = 14000
$ws.Range("J40").Value = 14068.571
$ws.Range("K40").Value = 14000
$ws.Range("L40").Value = 14068.571
$ws.Range("M40").Value = -13849
$ws.Range("N40").Value = -14370.571
$ws.Range("H113").Value = 72955.94
$ws.Range("I113").Value = 141848
$ws.Range("J113").Value = 4063.875
$ws.Range("K113").Value = 141848
$ws.Range("L113").Value = 4063.875
$ws.Range("M113").Value = -139678
$ws.Range("N113").Value = -8403.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 135889.73
$ws.Range("I22").Value = 3638.111
$ws.Range("J22").Value = 334267.16
$ws.Range("K22").Value = 3638.111
$ws.Range("L22").Value = 334267.16
$ws.Range("M22").Value = -3343.111
$ws.Range("N22").Value = -334857.16
$ws.Range("H27").Value = 135889.73
$ws.Range("I27").Value = 3638.111
$ws.Range("J27").Value = 334267.16
$ws.Range("K27").Value = 3638.111
$ws.Range("L27").Value = 334267.16
$ws.Range("M27").Value = -3531.111
$ws.Range("N27").Value = -334481.16
$ws.Range("H46").Value = 2097.0688
$ws.Range("I46").Value = 1279.5834
$ws.Range("J46").Value = 2674.1177
$ws.Range("K46").Value = 1279.5834
$ws.Range("L46").Value = 2674.1177
$ws.Range("M46").Value = -1091.5834
$ws.Range("N46").Value = -3050.1177
$ws.Range("H112").Value = 54387
$ws.Range("J112").Value = 54387
$ws.Range("L112").Value = 54387
$ws.Range("N112").Value = -57341
$ws.Range("H132").Value = 3462.7666
$ws.Range("I132").Value = 2562.1
$ws.Range("J132").Value = 5264.1
$ws.Range("K132").Value = 7686.299999999999
$ws.Range("L132").Value = 15792.3
$ws.Range("M132").Value = -5156.299999999999
$ws.Range("N132").Value = -20852.3
$ws.Range("H136").Value = 7767.857
$ws.Range("I136").Value = 10969.5
$ws.Range("K136").Value = 32908.5
$ws.Range("M136").Value = -30358.5
$ws.Range("H139").Value = 99138.86
$ws.Range("J139").Value = 99138.86
$ws.Range("L139").Value = 99138.86
$ws.Range("N139").Value = -109418.86

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 5850919.5
$ws.Range("I96").Value = 5000
$ws.Range("J96").Value = 8773879
$ws.Range("K96").Value = 5000
$ws.Range("L96").Value = 8773879
$ws.Range("M96").Value = -3627
$ws.Range("N96").Value = -8776625
$ws.Range("H132").Value = 2717.6428
$ws.Range("I132").Value = 2420.5833
$ws.Range("K132").Value = 7261.749899999999
$ws.Range("M132").Value = -4731.749899999999
